# feature/ADMINDASH-552: Add new reason in import attempt
#
# Sheet2 holds the master list of "reason" codes used by the dropdown
# validations on Sheet1. This change removes the obsolete
# DRIVER_ARRIVED_TOO_LATE reason and introduces two new reasons,
# CONSIGNEE_UNKNOWN and NO_CONSIGNEE, while re-ordering/grouping the
# remaining reasons.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")

# --- Rewrite Sheet2's reason list (now 29 rows, was 28) ---
$ws2.Cells.Item(1, 1).Value = "BAD_ADDRESS"
$ws2.Cells.Item(2, 1).Value = "CONSIGNEE_NOT_AROUND"
$ws2.Cells.Item(3, 1).Value = "CONSIGNEE_REFUSED_TO_ACCEPT"
$ws2.Cells.Item(4, 1).Value = "CONSIGNEE_CANNOT_BE_CONTACTED"
$ws2.Cells.Item(5, 1).Value = "CONSIGNEE_DOES_NOT_HAVE_ENOUGH_CASH"
$ws2.Cells.Item(6, 1).Value = "STUFF_OR_BOX_IS_BROKEN"
$ws2.Cells.Item(7, 1).Value = "CONSIGNEE_WANT_RESCHEDULE"
$ws2.Cells.Item(8, 1).Value = "NATURAL_DISASTER"
$ws2.Cells.Item(9, 1).Value = "CONSIGNEE_MOVE_OUT"
$ws2.Cells.Item(10, 1).Value = "CONSIGNEE_NOT_RECOGNIZED"
$ws2.Cells.Item(11, 1).Value = "CONSIGNEE_RARELY_IN_PLACE"
$ws2.Cells.Item(12, 1).Value = "OFFICE_CLOSED"
$ws2.Cells.Item(13, 1).Value = "LEAVE_HOLIDAY_SICK"
$ws2.Cells.Item(14, 1).Value = "OUT_OF_TOWN"
$ws2.Cells.Item(15, 1).Value = "CONSIGNEE_PASSED_AWAY"
$ws2.Cells.Item(16, 1).Value = "CONSIGNEE_RETIRED_RESIGNED"
$ws2.Cells.Item(17, 1).Value = "CONSIGNEE_DIFFICULT_TO_MEET"
$ws2.Cells.Item(18, 1).Value = "NEGATIVE_LOCATION"
$ws2.Cells.Item(19, 1).Value = "INCOMPLETE_ADDRESS"
$ws2.Cells.Item(20, 1).Value = "STUFF_DOES_NOT_MATCH_SPECIFICATION"
$ws2.Cells.Item(21, 1).Value = "EKYC_FAILED"
$ws2.Cells.Item(22, 1).Value = "PACKAGE_NOT_READY"
$ws2.Cells.Item(23, 1).Value = "PACKAGE_OVERSIZED"
$ws2.Cells.Item(24, 1).Value = "CONSIGNEE_UNKNOWN"
$ws2.Cells.Item(25, 1).Value = "NO_CONSIGNEE"
$ws2.Cells.Item(26, 1).Value = "COD_MISMATCH"
$ws2.Cells.Item(27, 1).Value = "MANUAL_PROCESS"
$ws2.Cells.Item(28, 1).Value = "OUT_OF_COVERAGE"
$ws2.Cells.Item(29, 1).Value = "CONSIGNEE_CHANGE_MIND"

# --- Refresh the dropdown validations so the bounded range keeps pace
#     with the longer reason list (Sheet2!$A$1:$A$19), and split the E2 /
#     B2:B1048576 combined rule into two independent ones. ---
$ws1.Range("E2").Validation.Delete()
$ws1.Range("B2:B1048576").Validation.Delete()
$ws1.Range("E3:E1048576").Validation.Delete()

$ws1.Range("E3:E1048576").Validation.Add(3, 1, 1, "Sheet2!`$A`$1:`$A`$19")
$ws1.Range("E2").Validation.Add(3, 1, 1, "Sheet2!`$A:`$A")
$ws1.Range("B2:B1048576").Validation.Add(3, 1, 1, "Sheet2!`$A:`$A")

# --- Restore the recorded selections / active sheet ---
$ws1.Range("E3").Select()
$ws2.Activate()
$ws2.Range("B20").Select()
